$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'isophonics_231'
$ws.Range("B2").Value = 'isophonics_124'
$ws.Range("C2").Value = 0.1087570621468927
$ws.Range("D2").Value = '[[''C'', ''G'', ''D'']]'
$ws.Range("E2").Value = '[[''Bb'', ''F'', ''C'']]'
$ws.Range("F2").Value = '[(9.236303, 12.162018)]'
$ws.Range("G2").Value = '[(16.749659, 21.811609)]'
$ws.Range("H2").Value = 'spotify:track:4F1AgKpuFRMLEgtPETVwZk'
$ws.Range("I2").Value = ""

# Row 3
$ws.Range("A3").Value = 'schubert-winterreise_129'
$ws.Range("B3").Value = 'schubert-winterreise_97'
$ws.Range("C3").Value = 0.1833333333333333
$ws.Range("D3").Value = '[[''C:7'', ''F:min'', ''C:maj''], [''F:maj'', ''A#:maj'', ''F:maj''], [''F:min'', ''C:maj'', ''F:min'']]'
$ws.Range("E3").Value = '[[''A:7'', ''D:min/A'', ''A:maj''], [''D:maj'', ''G:maj'', ''D:maj/F#''], [''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("F3").Value = '[(12.66, 17.1), (105.0, 108.24), (0.78, 5.26)]'
$ws.Range("G3").Value = '[(33.84, 36.48), (69.36, 73.2), (8.54, 14.08)]'
$ws.Range("H3").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I3").Value = ""

# Row 4
$ws.Range("A4").Value = 'jaah_7'
$ws.Range("B4").Value = 'schubert-winterreise_128'
$ws.Range("C4").Value = 0.03263106320431161
$ws.Range("D4").Value = '[[''Eb'', ''F:min7'', ''Bb:7''], [''F:min7'', ''Bb:7'', ''Eb''], [''Bb:7'', ''Eb'', ''Eb'']]'
$ws.Range("E4").Value = '[[''G:maj/D'', ''A:min7/C'', ''D:7''], [''A:min7/C'', ''D:7'', ''G:maj''], [''D:7'', ''G:maj'', ''G:maj'']]'
$ws.Range("F4").Value = '[(10.7, 14.24), (8.35, 10.7), (8.94, 11.88)]'
$ws.Range("G4").Value = '[(73.16, 77.72), (74.2, 81.18), (6.98, 13.38)]'
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 5
$ws.Range("A5").Value = 'isophonics_21'
$ws.Range("B5").Value = 'schubert-winterreise_161'
$ws.Range("C5").Value = 0.2554347826086957
$ws.Range("D5").Value = '[[''C'', ''G/3'', ''C'', ''G/3'', ''C'', ''G/3'']]'
$ws.Range("E5").Value = '[[''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'', ''C:maj'', ''G:maj'']]'
$ws.Range("F5").Value = '[(130.148, 132.444)]'
$ws.Range("G5").Value = '[(11.86, 15.88)]'
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("A6").Value = 'isophonics_232'
$ws.Range("B6").Value = 'isophonics_139'
$ws.Range("C6").Value = 0.2163742690058479
$ws.Range("D6").Value = '[[''C'', ''F'', ''G'', ''A'']]'
$ws.Range("E6").Value = '[[''G'', ''C'', ''D'', ''G'']]'
$ws.Range("F6").Value = '[(23.725555, 34.801473)]'
$ws.Range("G6").Value = '[(50.506553, 57.05458)]'
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = 'spotify:track:25yQPHgC35WNnnOUqFhgVR'

# Row 7
$ws.Range("A7").Value = 'isophonics_255'
$ws.Range("B7").Value = 'schubert-winterreise_6'
$ws.Range("C7").Value = 0.08684863523573201
$ws.Range("D7").Value = '[[''D:min'', ''G:min/5'', ''D:min'']]'
$ws.Range("E7").Value = '[[''B:min'', ''E:min/B'', ''B:min'']]'
$ws.Range("F7").Value = '[(27.789047, 31.109501)]'
$ws.Range("G7").Value = '[(27.36, 34.92)]'
$ws.Range("H7").Value = 'spotify:track:6rHh8urosEFRI67xVa6fzU'
$ws.Range("I7").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 8
$ws.Range("A8").Value = 'schubert-winterreise_157'
$ws.Range("B8").Value = 'jaah_1'
$ws.Range("C8").Value = 0.1669565217391304
$ws.Range("D8").Value = '[[''F:maj/A'', ''C:7'', ''F:maj'', ''F:maj/A''], [''F:maj'', ''F:maj/A'', ''C:7'', ''F:maj'']]'
$ws.Range("E8").Value = '[[''Eb'', ''Bb:7'', ''Eb'', ''Eb''], [''Eb'', ''Eb'', ''Bb:7'', ''Eb'']]'
$ws.Range("F8").Value = '[(16.2, 19.54), (17.94, 20.7)]'
$ws.Range("G8").Value = '[(4.83, 12.97), (2.77, 10.94)]'
$ws.Range("H8").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'
$ws.Range("I8").Value = ""

# Row 9
$ws.Range("A9").Value = 'isophonics_49'
$ws.Range("B9").Value = 'isophonics_220'
$ws.Range("C9").Value = 0.2476190476190476
$ws.Range("D9").Value = '[[''G'', ''G'', ''G'']]'
$ws.Range("E9").Value = '[[''D'', ''D'', ''D/7'']]'
$ws.Range("F9").Value = '[(20.940758, 23.936132)]'
$ws.Range("G9").Value = '[(0.325509, 3.175895)]'
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("A10").Value = 'schubert-winterreise_37'
$ws.Range("B10").Value = 'schubert-winterreise_78'
$ws.Range("C10").Value = 0.3095238095238095
$ws.Range("D10").Value = '[[''F:min/C'', ''C'', ''F:min/C'', ''C'', ''F:min'', ''C'']]'
$ws.Range("E10").Value = '[[''D:min'', ''A:maj'', ''D:min'', ''A:maj'', ''D:min'', ''A:maj'']]'
$ws.Range("F10").Value = '[(45.58, 53.44)]'
$ws.Range("G10").Value = '[(7.4, 18.96)]'
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = ""

# Row 11
$ws.Range("A11").Value = 'isophonics_157'
$ws.Range("B11").Value = 'isophonics_273'
$ws.Range("C11").Value = 0.1704545454545455
$ws.Range("D11").Value = '[[''E'', ''B'', ''E'']]'
$ws.Range("E11").Value = '[[''C'', ''G'', ''C'']]'
$ws.Range("F11").Value = '[(12.921927, 21.908049)]'
$ws.Range("G11").Value = '[(20.801, 27.165)]'
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

# Row 12
$ws.Range("A12").Value = 'schubert-winterreise_200'
$ws.Range("B12").Value = 'schubert-winterreise_25'
$ws.Range("C12").Value = 0.2589285714285714
$ws.Range("D12").Value = '[[''E:min'', ''B:maj'', ''E:min'']]'
$ws.Range("E12").Value = '[[''F:min/C'', ''C'', ''F:min/C'']]'
$ws.Range("F12").Value = '[(81.28, 89.42)]'
$ws.Range("G12").Value = '[(42.52, 46.36)]'
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

# Row 13
$ws.Range("A13").Value = 'schubert-winterreise_6'
$ws.Range("B13").Value = 'schubert-winterreise_192'
$ws.Range("C13").Value = 0.1098901098901099
$ws.Range("D13").Value = '[[''F#:7/B'', ''B:min'', ''E:min/B''], [''B:min'', ''F#:maj'', ''B:min'']]'
$ws.Range("E13").Value = '[[''C:7'', ''F:min'', ''A#:min/C#''], [''F:min/C'', ''C'', ''F:min/C'']]'
$ws.Range("F13").Value = '[(24.76, 32.42), (13.98, 21.5)]'
$ws.Range("G13").Value = '[(25.62, 33.76), (47.68, 51.98)]'
$ws.Range("H13").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Range("I13").Value = ""

# Row 14
$ws.Range("A14").Value = 'schubert-winterreise_61'
$ws.Range("B14").Value = 'isophonics_155'
$ws.Range("C14").Value = 0.06593406593406594
$ws.Range("D14").Value = '[[''C:maj/E'', ''G:maj/D'', ''A:min7/C''], [''G:min/D'', ''D:7'', ''G:min''], [''G:maj'', ''C:maj/G'', ''G:maj'']]'
$ws.Range("E14").Value = '[[''Db'', ''Ab'', ''Db:maj6/2''], [''F:min'', ''C:7'', ''F:min''], [''Ab'', ''Db/5'', ''Ab'']]'
$ws.Range("F14").Value = '[(66.02, 68.98), (42.44, 44.36), (18.12, 24.54)]'
$ws.Range("G14").Value = '[(63.858, 80.628), (17.833, 21.409), (0.243, 9.12)]'
$ws.Range("H14").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I14").Value = ""

# Row 15
$ws.Range("A15").Value = 'schubert-winterreise_121'
$ws.Range("B15").Value = 'jaah_69'
$ws.Range("C15").Value = 0.06554054054054054
$ws.Range("D15").Value = '[[''C'', ''C/G'', ''G:7'']]'
$ws.Range("E15").Value = '[[''Eb'', ''Eb'', ''Bb:7'']]'
$ws.Range("F15").Value = '[(282.12, 284.72)]'
$ws.Range("G15").Value = '[(33.48, 42.01)]'
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_128'
$ws.Range("B16").Value = 'isophonics_168'
$ws.Range("C16").Value = 0.06474358974358974
$ws.Range("D16").Value = '[[''G:maj'', ''E:min'', ''A:min/C'']]'
$ws.Range("E16").Value = '[[''G'', ''E:min'', ''A:min'']]'
$ws.Range("F16").Value = '[(10.96, 17.02)]'
$ws.Range("G16").Value = '[(37.99288, 42.009932)]'
$ws.Range("H16").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I16").Value = ""

# Row 17
$ws.Range("A17").Value = 'schubert-winterreise_199'
$ws.Range("B17").Value = 'schubert-winterreise_65'
$ws.Range("C17").Value = 0.3142857142857143
$ws.Range("D17").Value = '[[''G:maj/D'', ''D:7'', ''G:maj''], [''G:maj/B'', ''C:maj'', ''G:maj/D'']]'
$ws.Range("E17").Value = '[[''D#:maj'', ''A#:7'', ''D#:maj''], [''D#:maj'', ''G#:maj'', ''D#:maj'']]'
$ws.Range("F17").Value = '[(63.84, 65.24), (62.2, 64.4)]'
$ws.Range("G17").Value = '[(43.54, 60.0), (127.6, 132.4)]'
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
